# Updated symbol list on Mon Jan 23 13:46:10 UTC 2023 with GitHub Actions
# Apply updated Price (column D) and Volume(1h) (column E) values for the crypto table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.43"
$ws.Range("E2").Value = "'0.69%"
$ws.Range("D3").Value = "'35.92"
$ws.Range("E3").Value = "'-3.03%"
$ws.Range("D4").Value = "'5.077"
$ws.Range("E4").Value = "'1.45%"
$ws.Range("D5").Value = "'0.07928"
$ws.Range("E5").Value = "'0.81%"
$ws.Range("D6").Value = "'2.138"
$ws.Range("E6").Value = "'-3.23%"
$ws.Range("D7").Value = "'7.961"
$ws.Range("E7").Value = "'-0.43%"
$ws.Range("E8").Value = "'3.03%"
$ws.Range("D9").Value = "'0.9236"
$ws.Range("E9").Value = "'0.27%"
$ws.Range("D10").Value = "'0.09757"
$ws.Range("E10").Value = "'2.21%"
$ws.Range("D11").Value = "'0.1854"
$ws.Range("E11").Value = "'-1.58%"
$ws.Range("D12").Value = "'0.08633"
$ws.Range("E12").Value = "'0.62%"
$ws.Range("D13").Value = "'0.03567"
$ws.Range("E13").Value = "'-0.91%"
$ws.Range("D14").Value = "'0.09942"
$ws.Range("E14").Value = "'-0.29%"
$ws.Range("D15").Value = "'0.001444"
$ws.Range("E15").Value = "'-2.16%"
$ws.Range("D16").Value = "'0.005723"
$ws.Range("E16").Value = "'0.74%"
$ws.Range("D17").Value = "'3.463"
$ws.Range("E17").Value = "'0.08%"
$ws.Range("E18").Value = "'22.24%"
$ws.Range("E20").Value = "'2.52%"
$ws.Range("D21").Value = "'5.176"
$ws.Range("E21").Value = "'8.51%"
$ws.Range("D22").Value = "'0.2213"
$ws.Range("E22").Value = "'0.62%"
$ws.Range("D23").Value = "'0.04554"
$ws.Range("E23").Value = "'-0.78%"
$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'0.58%"
$ws.Range("D25").Value = "'0.004885"
$ws.Range("E25").Value = "'9.30%"
$ws.Range("D26").Value = "'0.0001303"
$ws.Range("E26").Value = "'-6.86%"
$ws.Range("D27").Value = "'0.0004761"
$ws.Range("E27").Value = "'0.21%"
$ws.Range("D39").Value = "'0.01858"
$ws.Range("E39").Value = "'1.67%"
$ws.Range("D40").Value = "'0.04743"
$ws.Range("E40").Value = "'0.38%"
$ws.Range("D41").Value = "'0.007920"
$ws.Range("E41").Value = "'-2.31%"
$ws.Range("D43").Value = "'0.007748"
$ws.Range("E43").Value = "'2.71%"
$ws.Range("D44").Value = "'0.002195"
$ws.Range("E44").Value = "'-0.64%"
$ws.Range("E45").Value = "'8.05%"
$ws.Range("D46").Value = "'0.00006278"
$ws.Range("E46").Value = "'1.90%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.40%"
$ws.Range("E48").Value = "'0.34%"
$ws.Range("D49").Value = "'49.96"
$ws.Range("E49").Value = "'366.54%"
$ws.Range("D50").Value = "'0.002005"
$ws.Range("E50").Value = "'-25.48%"
$ws.Range("D51").Value = "'0.00002105"
$ws.Range("E51").Value = "'0.40%"